$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Idea Planner")

# Update the "Done?" status for Phase 1.2 tasks (rows 13-16 in the Tasks table)
$ws.Range("E13").Value = "Yes"
$ws.Range("E14").Value = "Pending"
$ws.Range("E15").Value = "Pending"
$ws.Range("E16").Value = "Pending"

# Force recalculation so the dependent Status Icon / Indicator formulas update
$excel.Calculate()

# Update the active selection on the sheet (no frozen/top-left scroll cell anymore)
$ws.Range("J13").Select()
